# V2.0 Upgrade Checklist - finish the .net 3.5 async methods column (F)
# for the rows that were missing it, and move the on-screen selection to
# reflect the newly-worked area of the "Methods" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Methods")

# Rows that already had "Done" through column E but were missing the
# ".net 3.5 Async" (column F) mark.
$rows = @(38, 43, 44, 81, 82, 83, 84, 88)
foreach ($r in $rows) {
    $ws.Range("F$r").Value = "Done"
}

# Bring the freshly completed rows into view / match the saved selection.
$ws.Activate()
$ws.Range("F89").Select()
